$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the seven shared-string texts (gender/number agreement fixes).
#    These strings are reused by multiple cells across the sheet, so changing
#    the underlying cell text anywhere it appears updates every occurrence
#    that shares the same string.
$ws.Range("D10").Value = "SYSTEM exibe a listagem das Competencias (Portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("D11").Value = "SYSTEM destaca a Competencia (Portfolio) selecionada na listagem"
$ws.Range("B12").Value = "Lider de Pessoas clica na opcao 'Editar' para modificar a Competencia (Portfolio) selecionada"
$ws.Range("B28").Value = "Lider de Pessoas clica na opcao 'Excluir' para excluir a Competencia (Portfolio) selecionada"
$ws.Range("D39").Value = "SYSTEM exibe a listagem das Competencias (Portfolio) com a Competencia (Portfolio) excluida"
$ws.Range("D59").Value = "SYSTEM exibe a listagem das Competencias (Portfolio) cadastradas apenas para visualizacao com a opcao 'Ajuda'"
$ws.Range("D83").Value = "SYSTEM exibe a listagem das Competencias (Portfolio) sem a Competencia (Portfolio) excluida"

# 2. TC6 (rows 66-74) gains a new step: Lider de Pessoas indicates "Sim" (not
#    "Nao") and therefore must also fill out the competency-levels table
#    before saving, which triggers a new row between the existing rows 72
#    and 73.
$ws.Range("B72").Value = "Lider de Pessoas indica Sim no campo 'Niveis estao modificados para esta competencia'"

$ws.Rows.Item(73).Insert()

$ws.Range("A73").Value = 8
$ws.Range("B73").Value = "Lider de Pessoas preenche os dados na tabela do campo 'Niveis da Competencia' com os novos nome, valor e descricao dos niveis de competencia"
$ws.Range("C73").Value = ""
$ws.Range("D73").Value = "SYSTEM apresenta a tabela no campo 'Niveis da Competencia' preenchida corretamente"
$ws.Range("E73").Value = ""
$ws.Range("F73").Value = ""

$ws.Range("A73").Font.Name = $ws.Range("A72").Font.Name()
$ws.Range("A73").Font.Size = $ws.Range("A72").Font.Size()
$ws.Range("A73").Font.Bold = $ws.Range("A72").Font.Bold()
$ws.Range("A73").Interior.Color = $ws.Range("A72").Interior.Color()
$ws.Range("A73").HorizontalAlignment = $ws.Range("A72").HorizontalAlignment()
$ws.Range("A73").VerticalAlignment = $ws.Range("A72").VerticalAlignment()
$ws.Range("A73").WrapText = $ws.Range("A72").WrapText()

$ws.Range("A74").Value = 9

Write-Host "done"
